$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Q7)
$ws.Range("B9").Value = 0.4602904635456296
$ws.Range("C9").Value = 0.4602904635456296
$ws.Range("D9").Value = 0.2367358528406057
$ws.Range("E9").Value = 0.4865550871593121
$ws.Range("F9").Value = 0.1727490967023163
$ws.Range("G9").Value = 6

# Row 10 (Q8)
$ws.Range("B10").Value = 0.3684953970074996
$ws.Range("C10").Value = 0.3684953970074996
$ws.Range("D10").Value = 0.15310089540861
$ws.Range("E10").Value = 0.3912810951331663
$ws.Range("F10").Value = 0.1611460725222393
$ws.Range("G10").Value = 3

# Row 11 (Q9)
$ws.Range("B11").Value = 0.570669944985061
$ws.Range("C11").Value = 0.570669944985061
$ws.Range("D11").Value = 0.3256641861092525
$ws.Range("E11").Value = 0.570669944985061
$ws.Range("F11").ClearContents()
$ws.Range("G11").Value = 1
